$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 229  # H4 was 255.5
$ws.Cells.Item(4, 9).Value = 186.8  # I4 was 209
$ws.Cells.Item(4, 10).Value = 440  # J4 was 395
$ws.Cells.Item(4, 11).Value = 186.8  # K4 was 209
$ws.Cells.Item(4, 12).Value = 440  # L4 was 395
$ws.Cells.Item(4, 13).Value = -72.80000000000001  # M4 was -95
$ws.Cells.Item(4, 14).Value = -668  # N4 was -623
$ws.Cells.Item(19, 8).Value = 1361.8462  # H19 was 1328.0714
$ws.Cells.Item(19, 10).Value = 1323.6666  # J19 was 1261.5714
$ws.Cells.Item(19, 12).Value = 1323.6666  # L19 was 1261.5714
$ws.Cells.Item(19, 14).Value = -1673.6666  # N19 was -1611.5714
$ws.Cells.Item(43, 8).Value = 6640.4  # H43 was 7100.25
$ws.Cells.Item(43, 10).Value = 6640.4  # J43 was 7100.25
$ws.Cells.Item(43, 12).Value = 6640.4  # L43 was 7100.25
$ws.Cells.Item(43, 14).Value = -6778.4  # N43 was -7238.25
$ws.Cells.Item(64, 8).Value = 17064.578  # H64 was 18483.883
$ws.Cells.Item(64, 9).Value = 22569  # I64 was 26082.7
$ws.Cells.Item(64, 11).Value = 22569  # K64 was 26082.7
$ws.Cells.Item(64, 13).Value = -22321  # M64 was -25834.7
$ws.Cells.Item(67, 8).Value = 17064.578  # H67 was 18483.883
$ws.Cells.Item(67, 9).Value = 22569  # I67 was 26082.7
$ws.Cells.Item(67, 11).Value = 22569  # K67 was 26082.7
$ws.Cells.Item(67, 13).Value = -21711  # M67 was -25224.7
$ws.Cells.Item(70, 8).Value = 3511.0386  # H70 was 3332.4285
$ws.Cells.Item(70, 9).Value = 1862.2222  # I70 was 1662.1666
$ws.Cells.Item(70, 10).Value = 4383.9414  # J70 was 4585.125
$ws.Cells.Item(70, 11).Value = 5586.6666  # K70 was 4986.4998
$ws.Cells.Item(70, 12).Value = 13151.8242  # L70 was 13755.375
$ws.Cells.Item(70, 13).Value = -5316.6666  # M70 was -4716.4998
$ws.Cells.Item(70, 14).Value = -13691.8242  # N70 was -14295.375
$ws.Cells.Item(73, 8).Value = 3511.0386  # H73 was 3332.4285
$ws.Cells.Item(73, 9).Value = 1862.2222  # I73 was 1662.1666
$ws.Cells.Item(73, 10).Value = 4383.9414  # J73 was 4585.125
$ws.Cells.Item(73, 11).Value = 5586.6666  # K73 was 4986.4998
$ws.Cells.Item(73, 12).Value = 13151.8242  # L73 was 13755.375
$ws.Cells.Item(73, 13).Value = -4650.6666  # M73 was -4050.4998
$ws.Cells.Item(73, 14).Value = -15023.8242  # N73 was -15627.375
$ws.Cells.Item(129, 8).Value = 2353.5  # H129 was 2400
$ws.Cells.Item(129, 9).Value = 2353.5  # I129 was 2400
$ws.Cells.Item(129, 11).Value = 7060.5  # K129 was 7200
$ws.Cells.Item(129, 13).Value = -2060.5  # M129 was -2200
$ws.Cells.Item(138, 8).Value = 29444.107  # H138 was 26190.785
$ws.Cells.Item(138, 9).Value = 2101.4707  # I138 was 2136.1765
$ws.Cells.Item(138, 10).Value = 52685.35  # J138 was 42547.92
$ws.Cells.Item(138, 11).Value = 6304.4121  # K138 was 6408.529500000001
$ws.Cells.Item(138, 12).Value = 158056.05  # L138 was 127643.76
$ws.Cells.Item(138, 13).Value = -1164.4121  # M138 was -1268.529500000001
$ws.Cells.Item(138, 14).Value = -168336.05  # N138 was -137923.76

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 29500.244  # H32 was 30311.084
$ws.Cells.Item(32, 9).Value = 33953.562  # I32 was 35038.84
$ws.Cells.Item(32, 11).Value = 33953.562  # K32 was 35038.84
$ws.Cells.Item(32, 13).Value = -33666.562  # M32 was -34751.84
$ws.Cells.Item(63, 8).Value = 2753.476  # H63 was 2821.2
$ws.Cells.Item(63, 9).Value = 2430.8823  # I63 was 2495.375
$ws.Cells.Item(63, 11).Value = 2430.8823  # K63 was 2495.375
$ws.Cells.Item(63, 13).Value = -1744.8823  # M63 was -1809.375
$ws.Cells.Item(66, 8).Value = 2753.476  # H66 was 2821.2
$ws.Cells.Item(66, 9).Value = 2430.8823  # I66 was 2495.375
$ws.Cells.Item(66, 11).Value = 12154.4115  # K66 was 12476.875
$ws.Cells.Item(66, 13).Value = -8722.411500000002  # M66 was -9044.875
$ws.Cells.Item(132, 8).Value = 1622.7727  # H132 was 2156.3845
$ws.Cells.Item(132, 9).Value = 1194.7894  # I132 was 1503.3
$ws.Cells.Item(132, 11).Value = 3584.3682  # K132 was 4509.9
$ws.Cells.Item(132, 13).Value = -1054.3682  # M132 was -1979.9

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(31, 8).Value = 4800  # H31 was 6999
$ws.Cells.Item(31, 9).Value = 0  # I31 was 6999
$ws.Cells.Item(31, 10).Value = 4800  # J31 was 0
$ws.Cells.Item(31, 11).Value = 0  # K31 was 6999
$ws.Cells.Item(31, 13).ClearContents()  # M31 was -6747
$ws.Cells.Item(31, 14).Value = -5304  # N31 was None
$ws.Cells.Item(86, 8).Value = 1823.9166  # H86 was 1656.3334
$ws.Cells.Item(86, 9).Value = 1673.375  # I86 was 1764
$ws.Cells.Item(86, 10).Value = 2125  # J86 was 1562.125
$ws.Cells.Item(86, 11).Value = 1673.375  # K86 was 1764
$ws.Cells.Item(86, 12).Value = 2125  # L86 was 1562.125
$ws.Cells.Item(86, 13).Value = -550.375  # M86 was -641
$ws.Cells.Item(86, 14).Value = -4371  # N86 was -3808.125
$ws.Cells.Item(89, 8).Value = 1823.9166  # H89 was 1656.3334
$ws.Cells.Item(89, 9).Value = 1673.375  # I89 was 1764
$ws.Cells.Item(89, 10).Value = 2125  # J89 was 1562.125
$ws.Cells.Item(89, 11).Value = 8366.875  # K89 was 8820
$ws.Cells.Item(89, 12).Value = 10625  # L89 was 7810.625
$ws.Cells.Item(89, 13).Value = -2750.875  # M89 was -3204
$ws.Cells.Item(89, 14).Value = -21857  # N89 was -19042.625
$ws.Cells.Item(105, 8).Value = 1429.697  # H105 was 1518.0294
$ws.Cells.Item(105, 9).Value = 1080.85  # I105 was 1061.8636
$ws.Cells.Item(105, 10).Value = 1966.3846  # J105 was 2354.3333
$ws.Cells.Item(105, 11).Value = 1080.85  # K105 was 1061.8636
$ws.Cells.Item(105, 12).Value = 1966.3846  # L105 was 2354.3333
$ws.Cells.Item(105, 13).Value = 666.1500000000001  # M105 was 685.1364000000001
$ws.Cells.Item(105, 14).Value = -5460.3846  # N105 was -5848.3333
$ws.Cells.Item(134, 8).Value = 11614.679  # H134 was 13273.792
$ws.Cells.Item(134, 9).Value = 18371.334  # I134 was 22715
$ws.Cells.Item(134, 10).Value = 3818.5386  # J134 was 3832.5833
$ws.Cells.Item(134, 11).Value = 55114.00199999999  # K134 was 68145
$ws.Cells.Item(134, 12).Value = 11455.6158  # L134 was 11497.7499
$ws.Cells.Item(134, 13).Value = -52579.00199999999  # M134 was -65610
$ws.Cells.Item(134, 14).Value = -16525.6158  # N134 was -16567.7499

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6252588  # H31 was 6669418.5
$ws.Cells.Item(31, 9).Value = 9091815  # I31 was 11111998
$ws.Cells.Item(31, 10).Value = 6288  # J31 was 5548.3335
$ws.Cells.Item(31, 11).Value = 9091815  # K31 was 11111998
$ws.Cells.Item(31, 12).Value = 6288  # L31 was 5548.3335
$ws.Cells.Item(31, 13).Value = -9091520  # M31 was -11111703
$ws.Cells.Item(31, 14).Value = -6878  # N31 was -6138.3335
$ws.Cells.Item(34, 8).Value = 6252588  # H34 was 6669418.5
$ws.Cells.Item(34, 9).Value = 9091815  # I34 was 11111998
$ws.Cells.Item(34, 10).Value = 6288  # J34 was 5548.3335
$ws.Cells.Item(34, 11).Value = 9091815  # K34 was 11111998
$ws.Cells.Item(34, 12).Value = 6288  # L34 was 5548.3335
$ws.Cells.Item(34, 13).Value = -9091613  # M34 was -11111796
$ws.Cells.Item(34, 14).Value = -6692  # N34 was -5952.3335
$ws.Cells.Item(58, 8).Value = 1726.4117  # H58 was 1773.875
$ws.Cells.Item(58, 9).Value = 1726.4117  # I58 was 1773.875
$ws.Cells.Item(58, 11).Value = 1726.4117  # K58 was 1773.875
$ws.Cells.Item(58, 13).Value = -1523.4117  # M58 was -1570.875
$ws.Cells.Item(86, 8).Value = 34632.81  # H86 was 34709.04
$ws.Cells.Item(86, 9).Value = 57427.383  # I86 was 53959.855
$ws.Cells.Item(86, 10).Value = 11838.23  # J86 was 12249.75
$ws.Cells.Item(86, 11).Value = 57427.383  # K86 was 53959.855
$ws.Cells.Item(86, 12).Value = 11838.23  # L86 was 12249.75
$ws.Cells.Item(86, 13).Value = -56304.383  # M86 was -52836.855
$ws.Cells.Item(86, 14).Value = -14084.23  # N86 was -14495.75
$ws.Cells.Item(89, 8).Value = 34632.81  # H89 was 34709.04
$ws.Cells.Item(89, 9).Value = 57427.383  # I89 was 53959.855
$ws.Cells.Item(89, 10).Value = 11838.23  # J89 was 12249.75
$ws.Cells.Item(89, 11).Value = 287136.915  # K89 was 269799.275
$ws.Cells.Item(89, 12).Value = 59191.14999999999  # L89 was 61248.75
$ws.Cells.Item(89, 13).Value = -281520.915  # M89 was -264183.275
$ws.Cells.Item(89, 14).Value = -70423.14999999999  # N89 was -72480.75
$ws.Cells.Item(132, 8).Value = 102305.9  # H132 was 68503.53
$ws.Cells.Item(132, 9).Value = 144151.28  # I132 was 84462.75
$ws.Cells.Item(132, 11).Value = 432453.84  # K132 was 253388.25
$ws.Cells.Item(132, 13).Value = -429923.84  # M132 was -250858.25
$ws.Cells.Item(134, 8).Value = 2627.0557  # H134 was 2621.611
$ws.Cells.Item(134, 9).Value = 2480.3125  # I134 was 2474.1875
$ws.Cells.Item(134, 11).Value = 7440.9375  # K134 was 7422.5625
$ws.Cells.Item(134, 13).Value = -4905.9375  # M134 was -4887.5625
$ws.Cells.Item(136, 8).Value = 1726.4117  # H136 was 1773.875
$ws.Cells.Item(136, 9).Value = 1726.4117  # I136 was 1773.875
$ws.Cells.Item(136, 11).Value = 5179.2351  # K136 was 5321.625
$ws.Cells.Item(136, 13).Value = -2629.2351  # M136 was -2771.625

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 716.2857  # H92 was 689.25
$ws.Cells.Item(92, 10).Value = 413  # J92 was 456.5
$ws.Cells.Item(92, 12).Value = 1239  # L92 was 1369.5
$ws.Cells.Item(92, 14).Value = -3735  # N92 was -3865.5
$ws.Cells.Item(114, 8).Value = 1062  # H114 was 0
$ws.Cells.Item(114, 9).Value = 125  # I114 was 0
$ws.Cells.Item(114, 10).Value = 1999  # J114 was 0
$ws.Cells.Item(114, 11).Value = 375  # K114 was 0
$ws.Cells.Item(114, 12).Value = 5997  # L114 was 0
$ws.Cells.Item(114, 13).Value = 2879  # M114 was None
$ws.Cells.Item(114, 14).Value = -12505  # N114 was None
$ws.Cells.Item(137, 8).Value = 3439  # H137 was 3128.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1287.7826  # H22 was 1253.2916
$ws.Cells.Item(22, 9).Value = 1189.9231  # I22 was 1190.6923
$ws.Cells.Item(22, 10).Value = 1415  # J22 was 1327.2727
$ws.Cells.Item(22, 11).Value = 1189.9231  # K22 was 1190.6923
$ws.Cells.Item(22, 12).Value = 1415  # L22 was 1327.2727
$ws.Cells.Item(22, 13).Value = -894.9231  # M22 was -895.6922999999999
$ws.Cells.Item(22, 14).Value = -2005  # N22 was -1917.2727
$ws.Cells.Item(27, 8).Value = 1287.7826  # H27 was 1253.2916
$ws.Cells.Item(27, 9).Value = 1189.9231  # I27 was 1190.6923
$ws.Cells.Item(27, 10).Value = 1415  # J27 was 1327.2727
$ws.Cells.Item(27, 11).Value = 1189.9231  # K27 was 1190.6923
$ws.Cells.Item(27, 12).Value = 1415  # L27 was 1327.2727
$ws.Cells.Item(27, 13).Value = -1082.9231  # M27 was -1083.6923
$ws.Cells.Item(27, 14).Value = -1629  # N27 was -1541.2727
$ws.Cells.Item(55, 8).Value = 1366.4546  # H55 was 1369
$ws.Cells.Item(55, 9).Value = 409.66666  # I55 was 357.25
$ws.Cells.Item(55, 10).Value = 1725.25  # J55 was 1947.1428
$ws.Cells.Item(55, 11).Value = 409.66666  # K55 was 357.25
$ws.Cells.Item(55, 12).Value = 1725.25  # L55 was 1947.1428
$ws.Cells.Item(55, 13).Value = -236.66666  # M55 was -184.25
$ws.Cells.Item(55, 14).Value = -2071.25  # N55 was -2293.1428
$ws.Cells.Item(68, 8).Value = 5722.1113  # H68 was 5000
$ws.Cells.Item(68, 9).Value = 5749.5  # I68 was 5000
$ws.Cells.Item(68, 10).Value = 5714.2856  # J68 was 0
$ws.Cells.Item(68, 11).Value = 5749.5  # K68 was 5000
$ws.Cells.Item(68, 12).Value = 5714.2856  # L68 was 0
$ws.Cells.Item(68, 13).Value = -5000.5  # M68 was -4251
$ws.Cells.Item(68, 14).Value = -7212.2856  # N68 was None
$ws.Cells.Item(71, 8).Value = 5722.1113  # H71 was 5000
$ws.Cells.Item(71, 9).Value = 5749.5  # I71 was 5000
$ws.Cells.Item(71, 10).Value = 5714.2856  # J71 was 0
$ws.Cells.Item(71, 11).Value = 28747.5  # K71 was 25000
$ws.Cells.Item(71, 12).Value = 28571.428  # L71 was 0
$ws.Cells.Item(71, 13).Value = -25003.5  # M71 was -21256
$ws.Cells.Item(71, 14).Value = -36059.428  # N71 was None
$ws.Cells.Item(136, 8).Value = 3563.0688  # H136 was 3626.0356
$ws.Cells.Item(136, 10).Value = 4137.4165  # J136 was 4349.909
$ws.Cells.Item(136, 12).Value = 12412.2495  # L136 was 13049.727
$ws.Cells.Item(136, 14).Value = -17512.2495  # N136 was -18149.727

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 3914.1  # H126 was 3954
$ws.Cells.Item(126, 9).Value = 3204.4375  # I126 was 3254.3125
$ws.Cells.Item(126, 11).Value = 9613.3125  # K126 was 9762.9375
$ws.Cells.Item(126, 13).Value = -7143.3125  # M126 was -7292.9375
